$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '61.863.11'
$ws.Range("E2").Value = '  +0.54%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.421.64'
$ws.Range("E3").Value = '  +0.87%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.998'
$ws.Range("E4").Value = '  -0.26%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '407.61'
$ws.Range("E5").Value = '  +0.96%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '128.15'
$ws.Range("E6").Value = '  -1.54%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.629'
$ws.Range("E7").Value = '  +6.78%  '

$ws.Range("E8").Value = '  -0.13%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.732'
$ws.Range("E9").Value = '  +7.52%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.138'
$ws.Range("E10").Value = '  +7.30%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '42.54'
$ws.Range("E11").Value = '  +2.62%  '

$ws.Range("B12").Value = 'Polkadot'
$ws.Range("C12").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '9.13'
$ws.Range("E12").Value = '  +9.75%  '

$ws.Range("B13").Value = 'Chainlink'
$ws.Range("C13").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '21.51'
$ws.Range("E13").Value = '  +9.07%  '

$ws.Range("B14").Value = 'TRON'
$ws.Range("C14").Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.141'
$ws.Range("E14").Value = '  +0.18%  '

$ws.Range("B15").Value = 'WrappedliquidstakedEther2.0'
$ws.Range("C15").Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '3.952.75'
$ws.Range("E15").Value = '  +0.53%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.0000203'
$ws.Range("E16").Value = '  +42.61%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '3.436.48'
$ws.Range("E17").Value = '  +1.27%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '12.32'
$ws.Range("E18").Value = '  +6.18%  '

$ws.Range("E19").Value = '  +6.50%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '61.727.43'
$ws.Range("E20").Value = '  +0.39%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '446.38'
$ws.Range("E21").Value = '  +43.72%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '91.62'
$ws.Range("E22").Value = '  +10.42%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '3.19'
$ws.Range("E23").Value = '  +1.13%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '12.94'
$ws.Range("E24").Value = '  +2.19%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '3.24'
$ws.Range("E25").Value = '  +3.16%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '32.83'
$ws.Range("E26").Value = '  +11.87%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '8.67'
$ws.Range("E27").Value = '  +8.08%  '

$ws.Range("E28").Value = '  -0.45%  '

$ws.Range("E29").Value = '  -4.35%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '2.80'
$ws.Range("E30").Value = '  +3.34%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '11.95'
$ws.Range("E31").Value = '  +6.24%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.171'
$ws.Range("E32").Value = '  -0.16%  '

$ws.Range("E33").Value = '  +0.29%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '42.59'
$ws.Range("E34").Value = '  -2.70%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.00'
$ws.Range("E35").Value = '  -0.08%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.0496'
$ws.Range("E36").Value = '  +3.47%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '53.28'
$ws.Range("E37").Value = '  +4.10%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.997'
$ws.Range("E38").Value = '  -0.18%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '3.36'
$ws.Range("E39").Value = '  +0.88%  '

$ws.Range("E40").Value = '  +8.23%  '

$ws.Range("E41").Value = '  -1.12%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.317'
$ws.Range("E42").Value = '  +2.01%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '141.57'
$ws.Range("E43").Value = '  +1.47%  '

$ws.Range("B44").Value = 'WEMIXToken'
$ws.Range("C44").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '2.55'
$ws.Range("E44").Value = '  +15.42%  '

$ws.Range("B45").Value = 'NEARProtocol'
$ws.Range("C45").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '4.18'
$ws.Range("E45").Value = '  +6.90%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '1.98'
$ws.Range("E46").Value = '  +1.86%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '16.56'
$ws.Range("E47").Value = '  -0.41%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '22.23'
$ws.Range("E48").Value = '  +5.74%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '3.762.86'
$ws.Range("E49").Value = '  +0.74%  '

$ws.Range("B50").Value = 'Cronos'
$ws.Range("C50").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.138'
$ws.Range("E50").Value = '  +19.81%  '

$ws.Range("B51").Value = 'ThetaToken'
$ws.Range("C51").Value = 'https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '2.09'
$ws.Range("E51").Value = '  +10.77%  '
